$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agenda")
$lo = $ws.ListObjects.Item(1)

$ws.Range("B7:E7").Insert(-4121)  # xlShiftDown = -4121
$lo.Resize($ws.Range("B2:E37"))

# Copy entire row 8 (full row) into row 7
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(7).PasteSpecial(-4122)  # xlPasteFormats

Write-Host "Done"
